# Timesheet update: "added time for Jigyas"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 corresponds to Jigyas, week of 44487 (first week)
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("G4").Value = 2

# Row 10 corresponds to Jigyas, week of 44494 (second week)
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 2
